$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range so we know how far down/right the data goes.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Column A holds the "Beteckning" (case id), used both for:
#  - locating the last data row
#  - the new second HYPERLINK() argument (friendly display text)
# Column C holds the "Förändrad" date, which moves from 45184 -> 45186
# for every data row (row 1 is the header).
$colA = 1
$colC = 3
$newDate = 45186

for ($r = 2; $r -le $lastRow; $r++) {

    $idCell = $ws.Cells.Item($r, $colA)
    $idValue = $idCell.Value2
    if ([string]::IsNullOrEmpty($idValue)) {
        continue
    }

    # --- Update the "Förändrad" date in column C ---
    $cCell = $ws.Cells.Item($r, $colC)
    if (-not [string]::IsNullOrEmpty($cCell.Value2)) {
        $cCell.Value = $newDate
    }

    # --- Add the friendly-name second argument to every HYPERLINK() formula on this row ---
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if (-not $cell.HasFormula) {
            continue
        }

        $formula = $cell.Formula
        if ($formula -notlike '*HYPERLINK(*') {
            continue
        }

        # Skip cells that already carry a second argument (idempotent re-run safety).
        $openParen = $formula.IndexOf('(')
        $closeParen = $formula.LastIndexOf(')')
        if ($openParen -lt 0 -or $closeParen -lt 0 -or $closeParen -le $openParen) {
            continue
        }
        $inner = $formula.Substring($openParen + 1, $closeParen - $openParen - 1)
        if ($inner -match '","') {
            continue
        }

        $newFormula = $formula.Substring(0, $closeParen) + ', "' + $idValue + '"' + $formula.Substring($closeParen)
        $cell.Formula = $newFormula
    }
}
